# Update "想去人数" (want-to-go count) values in column F on the
# "展览" sheet and the "全部类型" sheet to match the newly scraped data.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 510
$wsExhibit.Range("F13").Value = 601
$wsExhibit.Range("F14").Value = 2579
$wsExhibit.Range("F15").Value = 125
$wsExhibit.Range("F16").Value = 3879
$wsExhibit.Range("F18").Value = 135
$wsExhibit.Range("F21").Value = 234
$wsExhibit.Range("F23").Value = 79
$wsExhibit.Range("F26").Value = 554
$wsExhibit.Range("F31").Value = 4314
$wsExhibit.Range("F34").Value = 320

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F8").Value = 510
$wsAll.Range("F16").Value = 601
$wsAll.Range("F18").Value = 2579
$wsAll.Range("F19").Value = 125
$wsAll.Range("F20").Value = 3879
$wsAll.Range("F22").Value = 135
$wsAll.Range("F25").Value = 234
$wsAll.Range("F28").Value = 79
$wsAll.Range("F31").Value = 554
$wsAll.Range("F36").Value = 4314
$wsAll.Range("F39").Value = 320
